$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New row 26 content
$ws.Range("A26").Value = "Get-Help * -Parameter computername | sort name | ft name, synopsis -auto -wrap"
$ws.Range("C26").Value = "Dit commando levert een lijst van alle`ncmdlets die native ondersteuning voor remoting hebben"

# Row height
$ws.Rows.Item(26).RowHeight = 45

# Merge A26:B26
$ws.Range("A26:B26").Merge()

# Styling to match existing pattern (horizontal left alignment added)
$ws.Range("A26:B26").HorizontalAlignment = -4131  # xlLeft

# Column B width adjustment
$ws.Columns.Item(2).ColumnWidth = 31.5703125

# Update selection to reflect new active cell A27
$ws.Range("A27").Select()
